$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row above the current first data row (row 7). This shifts
# the existing data rows (7-11) down to (8-12) and the footer rows (12-13)
# down to (13-14), adjusting merged cells automatically.
$ws.Rows("7:7").Insert()

# Copy the formatting of the row that used to be row 7 (now row 8) into the
# newly inserted blank row 7 so the new row matches the report's styling.
$ws.Range("A8:Q8").Copy()
$ws.Range("A7:Q7").PasteSpecial(-4122)

# Re-create the merged cells for the new row 7 (Insert() does not merge the
# freshly inserted blank row automatically).
$ws.Range("A7:B7").Merge()
$ws.Range("C7:G7").Merge()
$ws.Range("H7:K7").Merge()
$ws.Range("L7:M7").Merge()
$ws.Range("N7:O7").Merge()

# Restore/normalize the report's alternating row heights (PasteSpecial does
# not carry row height, and Insert() shifts the old heights down with the
# content instead of leaving the position-based pattern the report uses).
$ws.Rows("7:7").RowHeight = 25.5
$ws.Rows("8:8").RowHeight = 24.75
$ws.Rows("9:9").RowHeight = 25.5
$ws.Rows("10:10").RowHeight = 24.75
$ws.Rows("11:11").RowHeight = 25.5
$ws.Rows("12:12").RowHeight = 25.5
$ws.Rows("13:13").RowHeight = 24.75

# Populate the new row with the new item's data. Numeric-looking values are
# prefixed with a leading apostrophe so they are stored as text (matching
# the other rows in this column, which are shared-string text too).
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = "AMARYL 3 MG 30 TABS"
$ws.Range("H7").Value = "0:1"
$ws.Range("L7").Value = "'1"
$ws.Range("N7").Value = "'87.00"
$ws.Range("P7").Value = "'28.7100"
$ws.Range("Q7").Value = "0:1"

# Renumber the "م" (index) column for the rows that were pushed down.
$ws.Range("A8").Value = 2
$ws.Range("A9").Value = 3
$ws.Range("A10").Value = 4
$ws.Range("A11").Value = 5
$ws.Range("A12").Value = 6

# Update the total of the "سعر البيع" column (now row 13) to include the new
# row's value.
$ws.Range("P13").Value = 231.97

# Update the generated timestamp shown in the report footer (now row 14).
$ws.Range("A14").Value = "Sunday, 28 September, 2025 11:35 AM"
